# Updates the cryptos price/volume table with freshly scraped values.
# Cells whose new text looks like a plain decimal number (e.g. "4.09") are
# written through a NumberFormat="@" (Text) round-trip so Excel keeps the
# exact literal digits/trailing zeros instead of silently converting the
# cell to a numeric value; the cell's original style is restored right
# after so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '35.227.88'
$ws.Range("E2").Value2 = '  +0.23%  '
$ws.Range("D3").Value2 = '1.892.20'
$ws.Range("E3").Value2 = '  +2.12%  '
$ws.Range("E4").Value2 = '  -0.25%  '
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '242.70'
$c.Style = $origStyle
$ws.Range("E5").Value2 = '  +1.93%  '
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.654'
$c.Style = $origStyle
$ws.Range("E6").Value2 = '  +5.41%  '
$ws.Range("E7").Value2 = '  -0.31%  '
$c = $ws.Range("D8")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '41.22'
$c.Style = $origStyle
$ws.Range("E8").Value2 = '  -1.48%  '
$ws.Range("E9").Value2 = '  +6.35%  '
$c = $ws.Range("D10")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '49.90'
$c.Style = $origStyle
$ws.Range("E10").Value2 = '  +7.30%  '
$c = $ws.Range("D11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.0707'
$c.Style = $origStyle
$ws.Range("E11").Value2 = '  +2.24%  '
$c = $ws.Range("D12")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.0995'
$c.Style = $origStyle
$ws.Range("E12").Value2 = '  +0.66%  '
$ws.Range("D13").Value2 = '2.168.84'
$ws.Range("E13").Value2 = '  +2.19%  '
$c = $ws.Range("D14")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '11.92'
$c.Style = $origStyle
$ws.Range("E14").Value2 = '  +4.71%  '
$c = $ws.Range("D15")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.690'
$c.Style = $origStyle
$ws.Range("E15").Value2 = '  +2.38%  '
$ws.Range("D16").Value2 = '1.884.97'
$ws.Range("E16").Value2 = '  +1.45%  '
$ws.Range("E17").Value2 = '  +2.28%  '
$ws.Range("D18").Value2 = '35.214.27'
$ws.Range("E18").Value2 = '  +0.23%  '
$c = $ws.Range("D19")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '71.08'
$c.Style = $origStyle
$ws.Range("E19").Value2 = '  +1.66%  '
$ws.Range("E20").Value2 = '  +2.29%  '
$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '240.57'
$c.Style = $origStyle
$ws.Range("E21").Value2 = '  +0.08%  '
$c = $ws.Range("D22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '12.37'
$c.Style = $origStyle
$ws.Range("E22").Value2 = '  +1.67%  '
$c = $ws.Range("D23")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '4.72'
$c.Style = $origStyle
$ws.Range("E23").Value2 = '  +0.20%  '
$ws.Range("E24").Value2 = '  -0.37%  '
$c = $ws.Range("D25")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '2.42'
$c.Style = $origStyle
$ws.Range("E25").Value2 = '  +31.93%  '
$c = $ws.Range("D26")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '2.28'
$c.Style = $origStyle
$ws.Range("E26").Value2 = '  +0.54%  '
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '170.12'
$c.Style = $origStyle
$ws.Range("E27").Value2 = '  +0.64%  '
$c = $ws.Range("D28")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '8.37'
$c.Style = $origStyle
$ws.Range("E28").Value2 = '  +4.83%  '
$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '18.18'
$c.Style = $origStyle
$ws.Range("E29").Value2 = '  +3.37%  '
$ws.Range("E30").Value2 = '  +2.38%  '
$c = $ws.Range("D31")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '4.10'
$c.Style = $origStyle
$ws.Range("E31").Value2 = '  +3.13%  '
$ws.Range("B32").Value2 = 'Hedera'
$ws.Range("C32").Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D32")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.0559'
$c.Style = $origStyle
$ws.Range("E32").Value2 = '  +0.91%  '
$ws.Range("B33").Value2 = 'ImmutableX'
$ws.Range("C33").Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D33")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.937'
$c.Style = $origStyle
$ws.Range("E33").Value2 = '  +16.09%  '
$ws.Range("E34").Value2 = '  -0.18%  '
$c = $ws.Range("D35")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '4.09'
$c.Style = $origStyle
$ws.Range("E35").Value2 = '  +2.03%  '
$ws.Range("E36").Value2 = '  -1.39%  '
$c = $ws.Range("D37")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '2.01'
$c.Style = $origStyle
$ws.Range("E37").Value2 = '  +0.81%  '
$ws.Range("E38").Value2 = '  +1.14%  '
$ws.Range("B39").Value2 = 'VeChain'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D39")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '0.0207'
$c.Style = $origStyle
$ws.Range("E39").Value2 = '  +3.59%  '
$ws.Range("B40").Value2 = 'ARBITRUM'
$ws.Range("C40").Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D40")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '1.08'
$c.Style = $origStyle
$ws.Range("E40").Value2 = '  +1.73%  '
$ws.Range("E41").Value2 = '  +15.23%  '
$c = $ws.Range("D42")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '15.88'
$c.Style = $origStyle
$ws.Range("E42").Value2 = '  +6.42%  '
$c = $ws.Range("D43")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '88.87'
$c.Style = $origStyle
$ws.Range("E43").Value2 = '  -0.93%  '
$ws.Range("D44").Value2 = '1.336.14'
$ws.Range("E44").Value2 = '  -0.34%  '
$c = $ws.Range("D45")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '48.29'
$c.Style = $origStyle
$ws.Range("E45").Value2 = '  +40.11%  '
$ws.Range("E46").Value2 = '  +2.41%  '
$ws.Range("E47").Value2 = '  -1.00%  '
$ws.Range("E48").Value2 = '  +1.40%  '
$c = $ws.Range("D49")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '6.50'
$c.Style = $origStyle
$ws.Range("E49").Value2 = '  +0.94%  '
$ws.Range("D50").Value2 = '2.078.08'
$ws.Range("E50").Value2 = '  +1.93%  '
$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value2 = '11.32'
$c.Style = $origStyle
$ws.Range("E51").Value2 = '  -13.04%  '
